$wb = $excel.ActiveWorkbook

# --- Rename the existing sheet, then duplicate it to create the second phase ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Fase 1"

$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Fase 2"

# --- Trim "Fase 2" down to just the first defect block (rows 1-27) ---
$ws2.Range("A28:A73").EntireRow.Delete()

# --- Re-size a few columns on "Fase 2" (best effort - Excel snaps to pixel widths) ---
$ws2.Columns.Item(5).ColumnWidth = 22.17
$ws2.Columns.Item(7).ColumnWidth = 18.83
$ws2.Columns.Item(9).ColumnWidth = 19.85

# --- Update the header date (E6) ---
$ws2.Range("E6").Value = 45968
$ws2.Range("E6").NumberFormat = "d\-m\-yy"

# --- Defect block 1 (row 10-12): Def-001 ---
$ws2.Range("A10").Value = 45957
$ws2.Range("E10").Value = "Control Asignaciones"
$ws2.Range("G10").Value = "Iniciacion/Estrategia"
$ws2.Range("I10").Value = "Iniciacion/Estrategia"
$ws2.Range("M10").Value = "1 hora"
$ws2.Range("B11").Value = "Se completaron unas cosas en el control de asignaciones"

# --- Defect block 2 (row 15-17): Def-002 ---
$ws2.Range("A15").Value = 45957
$ws2.Range("E15").Value = "Acta de Iniciacion"
$ws2.Range("G15").Value = "Iniciacion/Estrategia"
$ws2.Range("I15").Value = "Iniciacion/Estrategia"
$ws2.Range("B16").Value = "Se le agrego lo que le faltaba"

# --- Defect block 3 (row 20-22): Def-003 ---
$ws2.Range("A20").Value = 45957
$ws2.Range("E20").Value = "PlanAdmin"
$ws2.Range("G20").Value = "Iniciacion/Estrategia"
$ws2.Range("I20").Value = "Iniciacion/Estrategia"
$ws2.Range("B21").Value = "Se agrego lo que faltaba"

# --- Defect block 4 (row 25-27): Def-004 ---
$ws2.Range("A25").Value = 45968
$ws2.Range("E25").Value = "Plan de calidad"
$ws2.Range("G25").Value = "Requerimiento/Planeacion"
$ws2.Range("I25").Value = "Requerimiento/Planeacion"
$ws2.Range("B26").Value = "Se hizo unos cambios"

Write-Output "done"
